$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find failed for: $old"
    }
}

function Replace-TextInRange($range, $old, $new) {
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find failed for: $old"
    }
}

# Title
Replace-Text "The Unseen Symphony: Exploring the Enigma of Quantum Entanglement" "The Significance of Civic Engagement in Shaping Our Democratic Future"

# Author
Replace-Text "Alex Thorn" "Hannah Mitchell"

# Email (the "io" -> "com" change must be scoped to the email paragraph,
# since "io" appears as a substring elsewhere in the document)
Replace-Text "thornscientist@scientist" "mitchellhan@privateconnect"
Replace-TextInRange $d.Paragraphs.Item(3).Range "io" "com"

# Body paragraph 1 (four sentences)
Replace-Text "In the vast symphony of existence, an enigmatic melody echoes through the scientific realm: quantum entanglement, a fascinating phenomenon where particles, separated by vast distances, remain intricately connected, defying the barriers of space and time" "Amidst the intricate tapestry of human societies, government stands as a cornerstone institution, orchestrating the collective efforts of individuals to achieve societal harmony and progress"

Replace-Text " This mind-bending dance of particles has captivated the imaginations of scientists, inspiring countless theories and experiments to unravel its secrets" " In the realm of governance, civic engagement emerges as a vital force, pulsating with the energy of informed and proactive citizens who actively participate in the decision-making processes that shape their communities"

Replace-Text " Quantum entanglement, a captivating paradox, challenges our understanding of reality, blurring the lines between what is real and what is merely perceived" " Civic engagement encompasses a diverse array of actions, ranging from casting a ballot in elections to volunteering in community projects, from advocating for policy changes to holding elected officials accountable"

Replace-Text " As we delve into this extraordinary phenomenon, we embark on a journey that promises to transform our understanding of the universe" " Through these endeavors, engaged citizens contribute to the vitality of their communities and safeguard the integrity of democratic principles"

# Body paragraph 1 continued after the double break
Replace-Text "From the subatomic ballet of fundamental particles to the vast cosmic web that binds celestial bodies, entanglement weaves its enigmatic tapestry across the universe" "In the intricate dance of civic engagement, citizens assume the mantle of active participants, rather than passive spectators, in the governance of their communities"

Replace-Text " Its implications extend far beyond the frontiers of scientific curiosity, probing the very foundations of reality" " They engage in informed dialogues with elected officials, advocating for policies that align with their values and aspirations"

Replace-Text " Envision two separated particles, each containing information" " By exercising their right to vote, they hold decision-makers accountable for their actions and shape the direction of public policy"

# Merge of three runs into one (delete middle runs, replace content)
Replace-Text " When one particle is observed, the state of its entangled counterpart instantaneously changes, regardless of the distance separating them. This eerie correlation, transcending the constraints of classical physics, suggests a profound non-local interconnectedness that defies our conventional notions of space, time, and causality" " Furthermore, engaged citizens actively participate in community projects, such as organizing neighborhood clean-ups or volunteering at local schools, fostering a sense of ownership and collective responsibility for their communities"

# Body paragraph 1 continued after the second double break
Replace-Text "Quantum entanglement, an ethereal chess game played by the universe itself, captivates us with its intricate moves" "The tapestry of civic engagement is woven from the threads of individual actions, each contributing to the vibrant fabric of a healthy democracy"

Replace-Text " Embarking on this intellectual quest, we seek to comprehend the orchestrator behind this unseen symphony, unravel the delicate mechanisms that govern the interplay of entangled particles" " When citizens actively participate in the governance of their communities, they create a virtuous cycle of accountability, transparency, and responsiveness"

Replace-Text " Is it an instantaneous exchange of information or an uncanny synchronization beyond the speed of light? Delving into this enigma, we may discover hidden dimensions, alternate realities, or even portals to other realms of existence" " Governments are more likely to enact policies that reflect the needs and desires of their constituents, leading to improved decision-making and better outcomes for all"

Replace-Text " The study of quantum entanglement holds the potential to reshape our understanding of reality, redefining the boundaries of the possible" " Civic engagement also fosters a sense of community and belonging, as citizens come together to address shared challenges and celebrate collective achievements"

# Summary heading (text unchanged, but touching it drops the lastRenderedPageBreak that sits in it)
Replace-Text "Summary" "Summary"

# Summary paragraph
Replace-Text "Quantum entanglement, an enigmatic phenomenon defying classical physics, captivates the scientific community, inspiring groundbreaking theories and experiments" "Civic engagement stands as a cornerstone of democratic governance, empowering citizens to actively participate in shaping their communities and ensuring the accountability of elected officials"

Replace-Text " Its implications extend beyond curiosity, challenging our notions of space, time, and causality" " Through informed dialogues, active participation in public affairs, and community involvement, engaged citizens contribute to the vibrancy and resilience of their communities"

Replace-Text " The exploration of quantum entanglement promises profound insights into the fundamental nature of the universe, potentially transforming our perception of reality" " Civic engagement promotes transparency, responsiveness, and accountability in government, leading to improved decision-making and better outcomes for all"

Replace-Text " This ethereal dance of particles invites us to question our understanding of existence and embark on a journey that may forever alter our comprehension of the cosmos" " By fostering a sense of community and collective responsibility, civic engagement strengthens the bonds that unite citizens and empowers them to shape a future that reflects their values and aspirations"

Write-Output "done"
